$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table on Sheet1 (A1:D60) is kept sorted descending by column D
# (h5_index). A new journal entry ("Journal of Democracy", h5_index 47)
# belongs right after the existing h5_index-52 row, i.e. at row 9 - so
# insert a fresh row there and push the rest of the table down by one.
$ws.Rows("9:9").Insert()

$ws.Range("A9").Value2 = "Journal of Democracy"
$ws.Range("B9").Value2 = "<a href='https://www.journalofdemocracy.org/about/submissions/'target='_blank'>Online Exclusive</a>"
$ws.Range("C9").Value2 = "2k words"
$ws.Range("D9").Value2 = 47

$ws.Range("I16").Select()
